$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13685
$ws1.Range("F4").Value = 126
$ws1.Range("F6").Value = 499
$ws1.Range("G6").Value = "不可售"
$ws1.Range("F9").Value = 13876
$ws1.Range("F10").Value = 14703
$ws1.Range("F18").Value = 17
$ws1.Range("F19").Value = 56
$ws1.Range("F21").Value = 1140
$ws1.Range("F24").Value = 5671
$ws1.Range("F25").Value = 941
$ws1.Range("F27").Value = 5388
$ws1.Range("F28").Value = 44
$ws1.Range("F30").Value = 232

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 4

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13685
$ws4.Range("F4").Value = 4
$ws4.Range("F5").Value = 126
$ws4.Range("F7").Value = 499
$ws4.Range("G7").Value = "不可售"
$ws4.Range("F10").Value = 13876
$ws4.Range("F11").Value = 14703
$ws4.Range("F19").Value = 17
$ws4.Range("F20").Value = 56
$ws4.Range("F22").Value = 1140
$ws4.Range("F25").Value = 5671
$ws4.Range("F26").Value = 941
$ws4.Range("F28").Value = 5388
$ws4.Range("F29").Value = 44
$ws4.Range("F31").Value = 232
